$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "62×95=5890" "70×66=4620"
Replace-Text "26×18=468" "42×53=2226"
Replace-Text "36×54=1944" "98×66=6468"
Replace-Text "93×99=9207" "32×42=1344"
Replace-Text "54×19=1026" "40×61=2440"
Replace-Text "31×46=1426" "63×69=4347"
Replace-Text "83×79=6557" "97×84=8148"
Replace-Text "77×45=3465" "15×14=210"
Replace-Text "15×17=255" "93×41=3813"
Replace-Text "59×68=4012" "12×72=864"
Replace-Text "23×40=920" "69×95=6555"
Replace-Text "98×94=9212" "69×92=6348"
Replace-Text "91×73=6643" "68×23=1564"
Replace-Text "50×17=850" "26×65=1690"
Replace-Text "92×39=3588" "21×42=882"
Replace-Text "65×31=2015" "73×97=7081"
Replace-Text "23×80=1840" "88×31=2728"
Replace-Text "26×82=2132" "91×62=5642"
Replace-Text "89×85=7565" "13×97=1261"
Replace-Text "52×34=1768" "97×32=3104"
Replace-Text "44×88=3872" "30×86=2580"
Replace-Text "68×79=5372" "45×71=3195"
Replace-Text "71×38=2698" "90×11=990"
Replace-Text "62×69=4278" "50×91=4550"
Replace-Text "98×86=8428" "73×79=5767"
